$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 75, shifting all existing rows (75-176) down by
# one — the data set gained one more weekly observation. This also pushes
# the former last row (176) down to 177, growing the used range to R177.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly record.
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 44467
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100114013
$ws.Range("G75").Value = "Zanahoria"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 160
$ws.Range("K75").Value = 6500
$ws.Range("L75").Value = 7000
$ws.Range("M75").Value = 6750
$ws.Range("N75").Value = "$/saco 20 kilos"
$ws.Range("O75").Value = "Provincia de Diguillín"
$ws.Range("P75").Value = 338
$ws.Range("Q75").Value = 20
$ws.Range("R75").Value = "Hortaliza"
